# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")

# Row 2 (H)
$wsOff.Range("B2").Value = 472
$wsOff.Range("C2").Value = 351
$wsOff.Range("D2").Value = 116
$wsOff.Range("E2").Value = 52

# Row 3 (R)
$wsOff.Range("B3").Value = 461
$wsOff.Range("C3").Value = 343
$wsOff.Range("D3").Value = 114
$wsOff.Range("E3").Value = 54
$wsOff.Range("F3").Value = 6

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")

# Row 2 (H)
$wsDef.Range("B2").Value = 449
$wsDef.Range("C2").Value = 311
$wsDef.Range("D2").Value = 96
$wsDef.Range("E2").Value = 41
$wsDef.Range("F2").Value = 6

# Row 3 (R)
$wsDef.Range("B3").Value = 448
$wsDef.Range("C3").Value = 313
$wsDef.Range("D3").Value = 108
$wsDef.Range("E3").Value = 54
